$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Phase 1: normalize formatting of the whole data block back to the
# plain bordered style (matches cellXfs index 1 in the target file) ---
$data = $ws.Range("A2:K4")
$data.Style = "Normal"
$data.Borders.Color = 0
$data.Borders.LineStyle = 1

# --- Phase 2: mark the code columns (AREA_CODE, BRANCH_CODE) as Text so
# that values with leading zeros ("01", "007", "012") are preserved ---
$ws.Range("A2:A4").NumberFormat = "@"
$ws.Range("C3:C4").NumberFormat = "@"

# --- Phase 3: write the real data, column by column, top to bottom, so
# that the shared-string table is rebuilt in the same order as the
# original authoring session ---
$ws.Range("A2").Value = "01"
$ws.Range("B2").Value = "Dhaka Area"
$ws.Range("D2").Value = "Ashkona Branch"

$ws.Range("A3").Value = "03"
$ws.Range("B3").Value = "Chattogram Area 1"
$ws.Range("C3").Value = "007"
$ws.Range("D3").Value = "Chandgaon Branch"

$ws.Range("A4").Value = "03"
$ws.Range("B4").Value = "Chattogram Area 1"
$ws.Range("C4").Value = "012"
$ws.Range("D4").Value = "Feni SME Branch"

$ws.Range("G2").Value = "BS0728"
$ws.Range("H2").Value = "Md. Shawkat Hossain"

$ws.Range("G3").Value = "BS1251"
$ws.Range("H3").Value = "Md. Jamshed Alam"

$ws.Range("G4").Value = "BS0901"
$ws.Range("H4").Value = "Mohammad Elias"

# C2 keeps a genuine numeric value (124) even though the column is
# formatted as text, so set the value before the number format.
$ws.Range("C2").Value = 124
$ws.Range("C2").NumberFormat = "@"

$ws.Range("J2").Value = 100
$ws.Range("J3").Value = 200
$ws.Range("J4").Value = 50

# --- Phase 4: tidy up the sheet view the same way the original author
# left it ---
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("K12").Select()
